$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated LR-pair stats (Natmi re-run "following Dr Hou advice"): the
# "Sending cluster"/"Target cluster" category labels were renumbered and
# every metric column (E..T) was recomputed, plus two brand-new cluster
# combinations (M2 as sender) were appended as rows 6-9.
$data = @(
    @("ECs",  "Gnai2", "Cxcr1", "ECs", 3, 1, 162.399297,       487.197891,   0.3910371682630009, 0.3910371682630009, 1, 0.3333333333333333, 2.161116666666667,  6.48335,  0.9584259542460944, 0.9584259542460944, 350.96382740165,   3158.67444661485, 0.3747801711381573,  0.3747801711381573),
    @("ECs",  "Gnai2", "Cxcr1", "M2",  3, 1, 162.399297,       487.197891,   0.3910371682630009, 0.3910371682630009, 3, 1,                  0.09374366666666666, 0.281231, 0.04157404575390553, 0.04157404575390552, 15.223905564869,   137.015150083821, 0.01625699712484365, 0.01625699712484365),
    @("FAPs", "Gnai2", "Cxcr1", "ECs", 3, 1, 65.41736466666667, 196.252094,  0.1575168212364948, 0.1575168212364948, 1, 0.3333333333333333, 2.161116666666667,  6.48335,  0.9584259542460944, 0.9584259542460944, 141.3745570705445, 1272.3710136349, 0.150968209703399,   0.150968209703399),
    @("FAPs", "Gnai2", "Cxcr1", "M2",  3, 1, 65.41736466666667, 196.252094,  0.1575168212364948, 0.1575168212364948, 3, 1,                  0.09374366666666666, 0.281231, 0.04157404575390553, 0.04157404575390552, 6.132463627523777, 55.19217264771399, 0.006548611533095793, 0.006548611533095793),
    @("M2",   "Gnai2", "Cxcr1", "ECs", 3, 1, 126.3069433333333, 378.92083,   0.3041313008456065, 0.3041313008456065, 1, 0.3333333333333333, 2.161116666666667,  6.48335,  0.9584259542460944, 0.9584259542460944, 272.9640403533889, 2456.6763631805,  0.2914873322290564, 0.2914873322290565),
    @("M2",   "Gnai2", "Cxcr1", "M2",  3, 1, 126.3069433333333, 378.92083,   0.3041313008456065, 0.3041313008456065, 3, 1,                  0.09374366666666666, 0.281231, 0.04157404575390553, 0.04157404575390552, 11.84047599352555, 106.56428394173,  0.01264396861655005, 0.01264396861655005),
    @("sCs",  "Gnai2", "Cxcr1", "ECs", 3, 1, 61.180387,         183.541161,  0.1473147096548978, 0.1473147096548978, 1, 0.3333333333333333, 2.161116666666667,  6.48335,  0.9584259542460944, 0.9584259542460944, 132.2179540188166, 1189.96158616935, 0.1411902411754818,  0.1411902411754818),
    @("sCs",  "Gnai2", "Cxcr1", "M2",  3, 1, 61.180387,         183.541161,  0.1473147096548978, 0.1473147096548978, 3, 1,                  0.09374366666666666, 0.281231, 0.04157404575390553, 0.04157404575390552, 5.735273805465666, 51.61746424919099, 0.006124468479416029, 0.006124468479416029)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}
